$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the two newly-identified market-closed holiday rows (row 7 = good Friday,
# row 13 = washington's bday) with -1 in the vol column, same convention used
# elsewhere in the sheet ("Input a -1 into vol to signify market closed day").
$ws.Range("B7").Value = -1
$ws.Range("B13").Value = -1

# Add the "Link to NYSE calendar" hyperlink + shared string under the existing
# "HOW TO USE" links list (G10:G12), matching their style/format.
$ws.Hyperlinks.Add($ws.Range("G13"), "https://www.nyse.com/markets/hours-calendars", "", "", "Link to NYSE calendar")
$ws.Range("G13").Style = $ws.Range("G12").Style

# Move the active selection to F11 to match the saved view state.
$ws.Range("F11").Select()
